$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "NCT(2.8169450852939946, 1.6038962399887462, -1.0900662187919659, 2.605752754704329)"
$ws.Range("C2").Value = "NIG(1.623920280084042, 1.3418022223479718, 3.220262939470181, 5.456774766805287)"
$ws.Range("D2").Value = "JSU(-0.7369152381565225, 1.0571326971150645, 0.6725926764094313, 3.000214114320405)"
$ws.Range("E2").Value = "NIG(0.779430990532496, 0.5057119149940451, 5.546002306028655, 4.694736544267396)"
